# Auto-generated edit script applying targeted cell updates
# per the Tiamat_Profits.xlsx diff (profit recalculation updates).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 7 (hunk 0)
$ws.Range("H7").Value = 20000
$ws.Range("J7").Value = 20000
$ws.Range("L7").Value = 20000
$ws.Range("N7").Value = -20224
# Row 10 (hunk 1)
$ws.Range("H10").Value = 44000
$ws.Range("I10").Value = 0
$ws.Range("J10").Value = 44000
$ws.Range("K10").Value = 0
$ws.Range("L10").Value = 44000
$ws.Range("M10").ClearContents()
$ws.Range("N10").Value = -44586
# Row 14 (hunk 2)
$ws.Range("H14").Value = 20000
$ws.Range("J14").Value = 20000
$ws.Range("L14").Value = 20000
$ws.Range("N14").Value = -20382
# Row 17 (hunk 3)
$ws.Range("H17").Value = 674019.9
$ws.Range("J17").Value = 674019.9
$ws.Range("L17").Value = 2022059.7
$ws.Range("N17").Value = -2022395.7
# Row 40 (hunk 4)
$ws.Range("H40").Value = 7131
$ws.Range("I40").Value = 9244.666999999999
$ws.Range("J40").Value = 790
$ws.Range("K40").Value = 9244.666999999999
$ws.Range("L40").Value = 790
$ws.Range("M40").Value = -9069.666999999999
$ws.Range("N40").Value = -1140
# Row 138 (hunk 5)
$ws.Range("H138").Value = 1828.45
$ws.Range("I138").Value = 676.7143
$ws.Range("J138").Value = 2276.3472
$ws.Range("K138").Value = 2030.1429
$ws.Range("L138").Value = 6829.0416
$ws.Range("M138").Value = 3109.8571
$ws.Range("N138").Value = -17109.0416

$ws = $wb.Worksheets.Item("BSM")
# Row 105 (hunk 6)
$ws.Range("H105").Value = 1501.2222
$ws.Range("I105").Value = 1333.3334
$ws.Range("J105").Value = 1585.1666
$ws.Range("K105").Value = 1333.3334
$ws.Range("L105").Value = 1585.1666
$ws.Range("M105").Value = 413.6666
$ws.Range("N105").Value = -5079.1666

$ws = $wb.Worksheets.Item("CRP")
# Row 10 (hunk 7)
$ws.Range("H10").Value = 70008
$ws.Range("I10").Value = 0
$ws.Range("K10").Value = 0
$ws.Range("M10").ClearContents()
# Row 31 (hunk 8)
$ws.Range("H31").Value = 5326.7637
$ws.Range("I31").Value = 4823
$ws.Range("J31").Value = 5778.4136
$ws.Range("K31").Value = 4823
$ws.Range("L31").Value = 5778.4136
$ws.Range("M31").Value = -4528
$ws.Range("N31").Value = -6368.4136
# Row 34 (hunk 9)
$ws.Range("H34").Value = 5326.7637
$ws.Range("I34").Value = 4823
$ws.Range("J34").Value = 5778.4136
$ws.Range("K34").Value = 4823
$ws.Range("L34").Value = 5778.4136
$ws.Range("M34").Value = -4621
$ws.Range("N34").Value = -6182.4136
# Row 60 (hunk 10)
$ws.Range("H60").Value = 5000
$ws.Range("I60").Value = 5000
$ws.Range("J60").Value = 0
$ws.Range("K60").Value = 5000
$ws.Range("L60").Value = 0
$ws.Range("M60").Value = -4489
$ws.Range("N60").ClearContents()
# Row 134 (hunk 11)
$ws.Range("H134").Value = 16667910
$ws.Range("I134").Value = 1245.6
$ws.Range("K134").Value = 3736.8
$ws.Range("M134").Value = -1201.8

$ws = $wb.Worksheets.Item("CUL")
# Row 68 (hunk 12)
$ws.Range("H68").Value = 1977.2727
$ws.Range("I68").Value = 2905.5898
$ws.Range("J68").Value = 1238.4082
$ws.Range("K68").Value = 8716.769400000001
$ws.Range("L68").Value = 3715.2246
$ws.Range("M68").Value = -7905.769400000001
$ws.Range("N68").Value = -5337.2246
# Row 69 (hunk 13)
$ws.Range("H69").Value = 18183878
$ws.Range("J69").Value = 18183878
$ws.Range("L69").Value = 54551634
$ws.Range("N69").Value = -54553256
# Row 71 (hunk 14)
$ws.Range("H71").Value = 1977.2727
$ws.Range("I71").Value = 2905.5898
$ws.Range("J71").Value = 1238.4082
$ws.Range("K71").Value = 26150.3082
$ws.Range("L71").Value = 11145.6738
$ws.Range("M71").Value = -22094.3082
$ws.Range("N71").Value = -19257.6738
# Row 72 (hunk 15)
$ws.Range("H72").Value = 18183878
$ws.Range("J72").Value = 18183878
$ws.Range("L72").Value = 163654902
$ws.Range("N72").Value = -163663014
# Row 104 (hunk 16)
$ws.Range("H104").Value = 4415.385
$ws.Range("I104").Value = 1000
$ws.Range("J104").Value = 4700
$ws.Range("K104").Value = 3000
$ws.Range("L104").Value = 14100
$ws.Range("M104").Value = -379
$ws.Range("N104").Value = -19342
# Row 121 (hunk 17)
$ws.Range("H121").Value = 810.9091
$ws.Range("I121").Value = 480
$ws.Range("J121").Value = 1000
$ws.Range("K121").Value = 1440
$ws.Range("L121").Value = 3000
$ws.Range("M121").Value = -130
$ws.Range("N121").Value = -5620

$ws = $wb.Worksheets.Item("GSM")
# Row 14 (hunk 18)
$ws.Range("H14").Value = 12667000
$ws.Range("I14").Value = 12667000
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 12667000
$ws.Range("L14").Value = 0
$ws.Range("M14").Value = -12666832
$ws.Range("N14").ClearContents()
# Row 21 (hunk 19)
$ws.Range("H21").Value = 2800
$ws.Range("I21").Value = 2800
$ws.Range("J21").Value = 0
$ws.Range("K21").Value = 2800
$ws.Range("L21").Value = 0
$ws.Range("M21").Value = -2627
$ws.Range("N21").ClearContents()
# Row 30 (hunk 20)
$ws.Range("H30").Value = 2800
$ws.Range("I30").Value = 2800
$ws.Range("J30").Value = 0
$ws.Range("K30").Value = 2800
$ws.Range("L30").Value = 0
$ws.Range("M30").Value = -2695
$ws.Range("N30").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
# Row 3 (hunk 21)
$ws.Range("H3").Value = 2538752.5
$ws.Range("J3").Value = 2898574.2
$ws.Range("L3").Value = 2898574.2
$ws.Range("N3").Value = -2898798.2
# Row 15 (hunk 22)
$ws.Range("H15").Value = 2538752.5
$ws.Range("J15").Value = 2898574.2
$ws.Range("L15").Value = 2898574.2
$ws.Range("N15").Value = -2898914.2
# Row 17 (hunk 23)
$ws.Range("H17").Value = 28847.715
$ws.Range("I17").Value = 14251.75
$ws.Range("J17").Value = 48309
$ws.Range("K17").Value = 14251.75
$ws.Range("L17").Value = 48309
$ws.Range("M17").Value = -14081.75
$ws.Range("N17").Value = -48649
# Row 18 (hunk 24)
$ws.Range("H18").Value = 11003
$ws.Range("J18").Value = 11003
$ws.Range("L18").Value = 11003
$ws.Range("N18").Value = -11347
# Row 22 (hunk 25)
$ws.Range("H22").Value = 613.2308
$ws.Range("I22").Value = 524.7273
$ws.Range("J22").Value = 1100
$ws.Range("K22").Value = 524.7273
$ws.Range("L22").Value = 1100
$ws.Range("M22").Value = -229.7273
$ws.Range("N22").Value = -1690
# Row 27 (hunk 26)
$ws.Range("H27").Value = 613.2308
$ws.Range("I27").Value = 524.7273
$ws.Range("J27").Value = 1100
$ws.Range("K27").Value = 524.7273
$ws.Range("L27").Value = 1100
$ws.Range("M27").Value = -417.7273
$ws.Range("N27").Value = -1314

$ws = $wb.Worksheets.Item("WVR")
# Row 133 (hunk 27)
$ws.Range("H133").Value = 43857.5
$ws.Range("J133").Value = 43857.5
$ws.Range("L133").Value = 43857.5
$ws.Range("N133").Value = -53977.5
# Row 136 (hunk 28)
$ws.Range("H136").Value = 401254.2
$ws.Range("I136").Value = 1346.4286
$ws.Range("J136").Value = 910227.75
$ws.Range("K136").Value = 4039.2858
$ws.Range("L136").Value = 2730683.25
$ws.Range("M136").Value = -1489.2858
$ws.Range("N136").Value = -2735783.25
# Row 141 (hunk 29)
$ws.Range("H141").Value = 65449
$ws.Range("J141").Value = 65449
$ws.Range("L141").Value = 65449
$ws.Range("N141").Value = -75809

Write-Host "Applied profit recalculation updates to ALC, BSM, CRP, CUL, GSM, LTW, WVR sheets."

